$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "63.320.97"
$c = $ws.Range("E2")
$c.NumberFormat = "@"
$c.Value = "  -1.24%  "
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "3.070.70"
$c = $ws.Range("E3")
$c.NumberFormat = "@"
$c.Value = "  -2.47%  "
$c = $ws.Range("E4")
$c.NumberFormat = "@"
$c.Value = "  +0.23%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "587.79"
$c = $ws.Range("E5")
$c.NumberFormat = "@"
$c.Value = "  -0.76%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "151.84"
$c = $ws.Range("E6")
$c.NumberFormat = "@"
$c.Value = "  +3.62%  "
$c = $ws.Range("E7")
$c.NumberFormat = "@"
$c.Value = "  +0.24%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.545"
$c = $ws.Range("E8")
$c.NumberFormat = "@"
$c.Value = "  +2.65%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "3.052.55"
$c = $ws.Range("E9")
$c.NumberFormat = "@"
$c.Value = "  -2.79%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.155"
$c = $ws.Range("E10")
$c.NumberFormat = "@"
$c.Value = "  -4.72%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "5.79"
$c = $ws.Range("E11")
$c.NumberFormat = "@"
$c.Value = "  -1.32%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.458"
$c = $ws.Range("E12")
$c.NumberFormat = "@"
$c.Value = "  -0.63%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "0.0000239"
$c = $ws.Range("E13")
$c.NumberFormat = "@"
$c.Value = "  -3.46%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "36.89"
$c = $ws.Range("E14")
$c.NumberFormat = "@"
$c.Value = "  -0.99%  "
$c = $ws.Range("E15")
$c.NumberFormat = "@"
$c.Value = "  -1.95%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "3.580.06"
$c = $ws.Range("E16")
$c.NumberFormat = "@"
$c.Value = "  -2.46%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "7.16"
$c = $ws.Range("E17")
$c.NumberFormat = "@"
$c.Value = "  -1.62%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "63.292.49"
$c = $ws.Range("E18")
$c.NumberFormat = "@"
$c.Value = "  -1.02%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "3.071.67"
$c = $ws.Range("E19")
$c.NumberFormat = "@"
$c.Value = "  -2.35%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "471.44"
$c = $ws.Range("E20")
$c.NumberFormat = "@"
$c.Value = "  +0.62%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "14.57"
$c = $ws.Range("E21")
$c.NumberFormat = "@"
$c.Value = "  +1.43%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "0.712"
$c = $ws.Range("E22")
$c.NumberFormat = "@"
$c.Value = "  -2.83%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "7.49"
$c = $ws.Range("E23")
$c.NumberFormat = "@"
$c.Value = "  -0.28%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "2.36"
$c = $ws.Range("E24")
$c.NumberFormat = "@"
$c.Value = "  +1.61%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "12.97"
$c = $ws.Range("E25")
$c.NumberFormat = "@"
$c.Value = "  -0.68%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "81.05"
$c = $ws.Range("E26")
$c.NumberFormat = "@"
$c.Value = "  -0.49%  "
$c = $ws.Range("B27")
$c.NumberFormat = "@"
$c.Value = "Dai"
$c = $ws.Range("C27")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "0.997"
$c = $ws.Range("E27")
$c.NumberFormat = "@"
$c.Value = "  -0.39%  "
$c = $ws.Range("B28")
$c.NumberFormat = "@"
$c.Value = "RenderToken"
$c = $ws.Range("C28")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "9.91"
$c = $ws.Range("E28")
$c.NumberFormat = "@"
$c.Value = "  +2.21%  "
$c = $ws.Range("B29")
$c.NumberFormat = "@"
$c.Value = "FirstDigitalUSD"
$c = $ws.Range("C29")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c = $ws.Range("E29")
$c.NumberFormat = "@"
$c.Value = "  -0.06%  "
$c = $ws.Range("B30")
$c.NumberFormat = "@"
$c.Value = "NEARProtocol"
$c = $ws.Range("C30")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "7.28"
$c = $ws.Range("E30")
$c.NumberFormat = "@"
$c.Value = "  -1.44%  "
$c = $ws.Range("B31")
$c.NumberFormat = "@"
$c.Value = "PancakeSwap"
$c = $ws.Range("C31")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "2.67"
$c = $ws.Range("E31")
$c.NumberFormat = "@"
$c.Value = "  -1.58%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "2.18"
$c = $ws.Range("E32")
$c.NumberFormat = "@"
$c.Value = "  -2.96%  "
$c = $ws.Range("E33")
$c.NumberFormat = "@"
$c.Value = "  +2.32%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "27.13"
$c = $ws.Range("E34")
$c.NumberFormat = "@"
$c.Value = "  -2.08%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "0.0₃0845"
$c = $ws.Range("E35")
$c.NumberFormat = "@"
$c.Value = "  +0.93%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "1.04"
$c = $ws.Range("E36")
$c.NumberFormat = "@"
$c.Value = "  -2.40%  "
$c = $ws.Range("B37")
$c.NumberFormat = "@"
$c.Value = "dogwifhat"
$c = $ws.Range("C37")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "3.33"
$c = $ws.Range("E37")
$c.NumberFormat = "@"
$c.Value = "  +2.62%  "
$c = $ws.Range("B38")
$c.NumberFormat = "@"
$c.Value = "Filecoin"
$c = $ws.Range("C38")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "6.07"
$c = $ws.Range("E38")
$c.NumberFormat = "@"
$c.Value = "  -1.76%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "2.20"
$c = $ws.Range("E39")
$c.NumberFormat = "@"
$c.Value = "  -4.85%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "9.25"
$c = $ws.Range("E40")
$c.NumberFormat = "@"
$c.Value = "  +0.34%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "50.24"
$c = $ws.Range("E41")
$c.NumberFormat = "@"
$c.Value = "  -2.16%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "444.87"
$c = $ws.Range("E42")
$c.NumberFormat = "@"
$c.Value = "  -3.64%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.283"
$c = $ws.Range("E43")
$c.NumberFormat = "@"
$c.Value = "  -3.22%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "0.0361"
$c = $ws.Range("E44")
$c.NumberFormat = "@"
$c.Value = "  -3.26%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "40.09"
$c = $ws.Range("E45")
$c.NumberFormat = "@"
$c.Value = "  +1.21%  "
$c = $ws.Range("E46")
$c.NumberFormat = "@"
$c.Value = "  +1.79%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "2.797.02"
$c = $ws.Range("E47")
$c.NumberFormat = "@"
$c.Value = "  -4.54%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "131.00"
$c = $ws.Range("E48")
$c.NumberFormat = "@"
$c.Value = "  +0.63%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "25.03"
$c = $ws.Range("E50")
$c.NumberFormat = "@"
$c.Value = "  +3.58%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "2.25"
$c = $ws.Range("E51")
$c.NumberFormat = "@"
$c.Value = "  -0.18%  "
